$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "theta_threshold_range" row (row 5) entirely, shifting rows up
$ws.Rows.Item(5).Delete()

# Update values for remaining rows
$ws.Range("B2").Value = 5.5
$ws.Range("C2").Value = 10

$ws.Range("B3").Value = 5.3
$ws.Range("C3").Value = 8.9

$ws.Range("B4").Value = 0.8
$ws.Range("C4").Value = 1.2

$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 15

# Update selection to match target (D10)
$ws.Range("D10").Select()
